# Gdnf-Gfra1.xlsx NATMI LR-pairs update ("Natmi following Dr Hou advice")
#
# A new target cluster, "ECs", is introduced. The row that used to describe
# the Gdnf -> Gfra1 signal onto "FAPs" becomes the "ECs" row (with refreshed
# NATMI statistics); "FAPs" then gets its own row with new numbers, and a
# third row is added for the original "sCs" target with numbers updated to
# match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 NATMI numeric metrics (columns E..T).
$rows = @(
    @(2, "sCs", "Gdnf", "Gfra1", "ECs",  3, 1, 1.995314333333334, 5.985943000000001, 1, 1, 2, 0.6666666666666666, 0.509196,           1.527588,  0.02558190413389134, 0.02558190413389134, 1.016006077276,     9.144054695484,      0.02558190413389134, 0.02558190413389134),
    @(3, "sCs", "Gdnf", "Gfra1", "FAPs", 3, 1, 1.995314333333334, 5.985943000000001, 1, 1, 3, 1,                  17.676258,          53.028774, 0.8880516296316739,  0.8880516296316739,  35.26969094709801,  317.427218523882,    0.8880516296316739,  0.8880516296316739),
    @(4, "sCs", "Gdnf", "Gfra1", "sCs",  3, 1, 1.995314333333334, 5.985943000000001, 1, 1, 3, 1,                  1.719084666666667, 5.157254,  0.0863664662344347,  0.0863664662344347,  3.430114275613556,  30.871028480522,     0.0863664662344347,  0.0863664662344347)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $r[$col]
    }
}
